$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("H9").Value = 67
$ws.Range("I9").Value = 67

$ws.Range("I9").Select()
